# Remove the stale "Ritu" test-user row from the "user" sheet.
# (Row 2 in the sheet - the TarakMehta row below it shifts up to row 2.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user")

# Delete the whole row; cells below shift up automatically.
$ws.Rows.Item(2).Delete()

# The hyperlinks on the old row 3 (now row 2) keep pointing at row 3 after
# the shift, so re-anchor them explicitly.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:nevixo9520@ociun.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:f9iupld30y@elatter.com")

# Active sheet/selection moved to the "user" sheet at E12.
$ws.Activate() | Out-Null
$ws.Range("E12").Select() | Out-Null
